$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D25").Value = "[바람돌이/딥러닝] GAN(2) - CGAN(Conditional Generative Adversarial Nets 논문 이론 및 리뷰"
$ws.Range("E25").Value = "https://blog.naver.com/winddori2002/222222304740"

$ws.Range("D37").Value = "[Paper Review] USAD: UnSupervised Anomaly Detection on Multivariate Time Series"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=1434&mod=document&pageid=1"

$ws.Range("D39").Value = "Face Alignment for Face Recognition in Python within OpenCV"
$ws.Range("E39").Value = "https://a292run.tistory.com/entry/Face-Alignment-for-Face-Recognition-in-Python-within-OpenCV-1"

$wb.Save()
